$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'26.742.30"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +0.46%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'1.641.04"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +0.01%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.36%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'217.85"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E6').Value = "'  +0.05%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'  +0.33%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = "'  +0.31%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('E9').Value = "'  -0.07%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'19.13"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  +0.16%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D12').Value = "'1.869.79"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  -0.06%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'1.639.40"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -0.08%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('E14').Value = "'  -0.31%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('E15').Value = "'  -0.34%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'64.70"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  -0.20%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'26.736.74"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +0.37%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('E18').Value = "'  -1.26%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'214.14"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  -0.42%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').Value = "'  +0.36%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('E21').Value = "'  +0.90%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('E22').Value = "'  +7.26%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'6.23"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  -0.45%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'9.28"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  -1.88%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'145.60"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  +0.27%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = "'  +0.32%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').Value = "'  -0.97%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = "'  +0.76%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'15.67"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  -0.06%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').Value = "'  -0.76%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').Value = "'  +1.43%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'3.40"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +1.24%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'3.00"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  +0.21%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'1.287.10"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  +0.82%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = "'  -0.07%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').Value = "'  +0.84%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = "'  -0.38%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = "'  +0.86%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'0.817"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  -0.51%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('E40').Value = "'  +0.35%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'0.806"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  -0.54%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').Value = "'  -1.16%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'5.29"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  -2.35%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'1.779.79"
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Value = "'61.07"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  +3.37%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('E46').Value = "'  +0.04%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'1.60"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  +0.03%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = "'  +0.46%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'7.60"
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Value = "'0.0965"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  +0.11%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'0.407"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +0.19%  "
$ws.Range('E51').Style = 'Normal'
